# C5-PowerPoint.pptx edit
#
# The only substantive, OOXML-visible change in the commit is on the table
# that lives on slide 6 ("SOURCES OF FINANCE"): its PowerPoint "quick style"
# (<a:tableStyleId>) is switched from the built-in "Medium Style 3 - Accent 1"
# GUID to the built-in "Medium Style 3" GUID. In the UI this is simply
# selecting a different style swatch from the Table Styles gallery while the
# table is selected.
#
# Table styles cannot be set through Table.Style (that property is
# read-only in this object model); PowerPoint exposes the mutation through
# Table.ApplyStyle("{GUID}") instead, so we look the table up and call that.

$p = $ppt.ActivePresentation

$newStyleId = "{19553013-BD97-4FE3-A08C-FAEE38C8DC51}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
